# PV_investment_electricity_cost.xlsx - "Updated for September 2020"
#
# The workbook tracks monthly PV (solar) production / electricity costs.
# Row 12 is the September row; it previously had only a few cells filled in
# (B12, N12, Q12, R12) with everything else at 0 / blank. This edit fills in
# the rest of September's figures (raw SolarEdge + E.On readings) and the
# per-row formulas that every other month row already has, mirroring the
# pattern used in row 11 (August).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- raw input data for row 12 (September) ------------------------------
$ws.Range("C12").Value = 963500
$ws.Range("D12").Value = 1021839
$ws.Range("E12").Value = 540904
$ws.Range("F12").Value = 422596
$ws.Range("G12").Value = 1444435

$ws.Range("I12").Value = 541
$ws.Range("J12").Value = 1021

$ws.Range("O12").Value = 383.75
$ws.Range("P12").Value = 49.69
$ws.Range("S12").Value = 48.24
$ws.Range("T12").Value = 2.92

$ws.Range("AD12").Value = 300

# --- formulas for row 12, matching the pattern already used in row 11 ---
$ws.Range("K12").Formula = "=(G12/1000)+I12-J12"
$ws.Range("L12").Formula = "=K12-I12"
$ws.Range("M12").Formula = "=L12/K12"
$ws.Range("U12").Formula = "=(K12*(P12+Q12+R12)/100)+N12+O12"
$ws.Range("V12").Formula = "=(I12*(P12+Q12+R12)/100)+N12+O12"
$ws.Range("W12").Formula = "=J12*(S12+T12)/100"
$ws.Range("X12").Formula = "=U12-V12"
$ws.Range("Y12").Formula = "=W12+X12"
$ws.Range("Z12").Formula = "=J12*0.6"
$ws.Range("AA12").Formula = "=Z12+Y12"
$ws.Range("AB12").Formula = "=AA12/(G12/1000)"
$ws.Range("AC12").Formula = "=(P12+Q12+R12)/100"

# --- leave the cursor where the author left it after the edit -----------
$ws.Range("A34").Select()
